$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "86.735.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.184.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.76%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "202.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -8.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "603.34"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -7.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.363"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -7.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.657"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +7.87%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.182.99"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.531"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -9.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.176"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -17.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.782.54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.08%  "

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.689.36"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.16%  "

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -9.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.181.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.35"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -9.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.91"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -8.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.61"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -9.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.41"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -13.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.02"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -8.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -9.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.40"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -11.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.364.67"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "73.23"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000128"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.167"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -16.81%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "533.93"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -11.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.27"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -12.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -14.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.26"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -22.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.53"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -10.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.133"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -9.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "21.73"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.81"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.96%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.373"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -11.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.86"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -15.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "145.94"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "172.38"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -9.78%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.15"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.93%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.126"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +10.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.25"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -13.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.95"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -10.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.588"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -11.65%  "

